$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 139. This pushes the existing rows
# 139-141 down to 141-143, preserving their original data intact.
$ws.Rows.Item(139).Insert()
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with the updated record.
$ws.Cells.Item(139, 1).Value2  = 9
$ws.Cells.Item(139, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(139, 3).Value2  = "Metropolitana"
$ws.Cells.Item(139, 4).Value2  = 44448
$ws.Cells.Item(139, 5).Value2  = 13
$ws.Cells.Item(139, 6).Value2  = 100112043
$ws.Cells.Item(139, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(139, 8).Value2  = "Sin especificar"
$ws.Cells.Item(139, 9).Value2  = "Primera"
$ws.Cells.Item(139, 10).Value2 = 79
$ws.Cells.Item(139, 11).Value2 = 17000
$ws.Cells.Item(139, 12).Value2 = 18000
$ws.Cells.Item(139, 13).Value2 = 17494
$ws.Cells.Item(139, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(139, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(139, 16).Value2 = 292
$ws.Cells.Item(139, 17).Value2 = 60
$ws.Cells.Item(139, 18).Value2 = "Hortaliza"

# Populate the new row 140 with the updated record.
$ws.Cells.Item(140, 1).Value2  = 9
$ws.Cells.Item(140, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(140, 3).Value2  = "Metropolitana"
$ws.Cells.Item(140, 4).Value2  = 44448
$ws.Cells.Item(140, 5).Value2  = 13
$ws.Cells.Item(140, 6).Value2  = 100112043
$ws.Cells.Item(140, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(140, 8).Value2  = "Sin especificar"
$ws.Cells.Item(140, 9).Value2  = "Segunda"
$ws.Cells.Item(140, 10).Value2 = 52
$ws.Cells.Item(140, 11).Value2 = 15000
$ws.Cells.Item(140, 12).Value2 = 16000
$ws.Cells.Item(140, 13).Value2 = 15500
$ws.Cells.Item(140, 14).Value2 = "$/caja 100 unidades"
$ws.Cells.Item(140, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(140, 16).Value2 = 155
$ws.Cells.Item(140, 17).Value2 = 100
$ws.Cells.Item(140, 18).Value2 = "Hortaliza"
